$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44203
$ws.Range("J2").Value = 20

# Row 3
$ws.Range("D3").Value = 44271
$ws.Range("J3").Value = 40
$ws.Range("K3").Value = 30000
$ws.Range("L3").Value = 30000
$ws.Range("M3").Value = 30000
$ws.Range("P3").Value = 1200

# Row 4
$ws.Range("D4").Value = 44264
$ws.Range("J4").Value = 30
$ws.Range("K4").Value = 29000
$ws.Range("L4").Value = 29000
$ws.Range("M4").Value = 29000
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 1160

# Row 5
$ws.Range("D5").Value = 44253
$ws.Range("J5").Value = 80

# Row 6
$ws.Range("D6").Value = 44222
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 40000
$ws.Range("L6").Value = 40000
$ws.Range("M6").Value = 40000
$ws.Range("P6").Value = 1600

# Row 7
$ws.Range("D7").Value = 44202
$ws.Range("J7").Value = 30

# Row 8
$ws.Range("D8").Value = 44246
$ws.Range("J8").Value = 50
$ws.Range("K8").Value = 31000
$ws.Range("L8").Value = 31000
$ws.Range("M8").Value = 31000
$ws.Range("P8").Value = 1240

# Row 9
$ws.Range("D9").Value = 44225
$ws.Range("J9").Value = 60
$ws.Range("K9").Value = 32000
$ws.Range("L9").Value = 32000
$ws.Range("M9").Value = 32000
$ws.Range("O9").Value = "Región del Maule"
$ws.Range("P9").Value = 1280

# Row 10
$ws.Range("D10").Value = 44232
$ws.Range("K10").Value = 40000
$ws.Range("L10").Value = 40000
$ws.Range("M10").Value = 40000
$ws.Range("P10").Value = 1600

# Row 11
$ws.Range("D11").Value = 44243
$ws.Range("K11").Value = 33000
$ws.Range("L11").Value = 33000
$ws.Range("M11").Value = 33000
$ws.Range("P11").Value = 1320

# Row 12
$ws.Range("D12").Value = 44250
$ws.Range("J12").Value = 70
$ws.Range("K12").Value = 30000
$ws.Range("L12").Value = 30000
$ws.Range("M12").Value = 30000
$ws.Range("P12").Value = 1200

# Row 13
$ws.Range("D13").Value = 44236
$ws.Range("J13").Value = 40
$ws.Range("K13").Value = 38000
$ws.Range("L13").Value = 38000
$ws.Range("M13").Value = 38000
$ws.Range("P13").Value = 1520

# Row 14
$ws.Range("D14").Value = 44201
$ws.Range("J14").Value = 60
$ws.Range("K14").Value = 30000
$ws.Range("L14").Value = 30000
$ws.Range("M14").Value = 30000
$ws.Range("P14").Value = 1200

# Row 15
$ws.Range("D15").Value = 44215
$ws.Range("K15").Value = 35000
$ws.Range("L15").Value = 35000
$ws.Range("M15").Value = 35000
$ws.Range("P15").Value = 1400

# Row 16
$ws.Range("D16").Value = 44204
$ws.Range("J16").Value = 50
$ws.Range("K16").Value = 37000
$ws.Range("L16").Value = 37000
$ws.Range("M16").Value = 37000
$ws.Range("P16").Value = 1480

# Row 17
$ws.Range("D17").Value = 44218
$ws.Range("J17").Value = 60
$ws.Range("K17").Value = 42000
$ws.Range("L17").Value = 42000
$ws.Range("M17").Value = 42000
$ws.Range("P17").Value = 1680

# Row 18
$ws.Range("D18").Value = 44239
$ws.Range("J18").Value = 60
$ws.Range("K18").Value = 35000
$ws.Range("L18").Value = 35000
$ws.Range("M18").Value = 35000
$ws.Range("P18").Value = 1400

# Row 19
$ws.Range("D19").Value = 44211
$ws.Range("J19").Value = 40
$ws.Range("K19").Value = 42000
$ws.Range("L19").Value = 42000
$ws.Range("M19").Value = 42000
$ws.Range("P19").Value = 1680
